$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- S09/G02 rows (81-83): style normalisation on deviation/remarks cells ---
# These cells already render with the default (general/bottom, no-wrap) look;
# make that explicit so the cell format is pinned rather than implicit.
foreach ($addr in @("F81","H81","F82","H82","F83","H83")) {
    $rng = $ws.Range($addr)
    $rng.VerticalAlignment = -4107   # xlVAlignBottom
    $rng.WrapText = $false
}

# --- Row 84: S09_G03_TB001 ---
$ws.Range("F84").Value2 = "require_admin now first honours a logged-in ADMIN user session, falling back to ST_ADMIN_USERNAME/ST_ADMIN_PASSWORD HTTP Basic only when configured."
$ws.Range("G84").Value2 = "implemented"
$ws.Range("H84").Value2 = "This keeps backwards compatibility for environments still using Basic auth while allowing the new auth system to be the primary way to access admin APIs."
$ws.Range("I84").Value2 = "In later sprints we may fully retire Basic auth and rely solely on user roles once all deployments have migrated."

# --- Row 85: S09_G03_TB002 ---
$ws.Range("F85").Value2 = "All admin-only routers (strategies, risk-settings, orders, positions, analytics, system-events, brokers) are wired through require_admin; Zerodha and webhook routes remain public."
$ws.Range("G85").Value2 = "implemented"
$ws.Range("H85").Value2 = "Broker configuration is again restricted to admins; normal users must log in as an ADMIN to edit secrets or risk settings."
$ws.Range("I85").Value2 = "Future work can add per-user views of orders/analytics separate from global admin data."

# --- Row 86: S09_G03_TB003 ---
$ws.Range("F86").Value2 = "Dev/test mode is handled by checking PYTEST_CURRENT_TEST in both get_settings and require_admin, keeping admin APIs open during pytest runs while enforcing auth in normal execution."
$ws.Range("G86").Value2 = "implemented"
$ws.Range("H86").Value2 = "This keeps the existing test suite stable while the UI begins to rely on session-based admin access in real usage."
$ws.Range("I86").Value2 = "Consider adding explicit auth-required tests once frontend flows are fully stabilised."
